$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "05/27/2015 Wed"
$ws.Range("C7").Value = "1:15pm-4:00pm"
$ws.Range("D7").Value = 2.75
$ws.Range("E7").Value = "Class Meeting, Team Meeting, Debugging"

$ws.Range("B8").Value = "05/28/2015 Thu"
$ws.Range("C8").Value = "4:45pm-"
$ws.Range("E8").Value = "SRS Document Revision,"

$ws.Range("E8").Select()
